$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4002304
$ws.Range("I137").Value = 5265031.5
$ws.Range("J137").Value = 3667.1667
$ws.Range("K137").Value = 15795094.5
$ws.Range("L137").Value = 11001.5001
$ws.Range("M137").Value = -15792544.5
$ws.Range("N137").Value = -16101.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3087.5
$ws.Range("I61").Value = 1042.5
$ws.Range("J61").Value = 8200
$ws.Range("K61").Value = 1042.5
$ws.Range("L61").Value = 8200
$ws.Range("M61").Value = -830.5
$ws.Range("N61").Value = -8624

$ws.Range("H63").Value = 2744.3635
$ws.Range("I63").Value = 2358.2856
$ws.Range("J63").Value = 3420
$ws.Range("K63").Value = 2358.2856
$ws.Range("L63").Value = 3420
$ws.Range("M63").Value = -1672.2856
$ws.Range("N63").Value = -4792

$ws.Range("H66").Value = 2744.3635
$ws.Range("I66").Value = 2358.2856
$ws.Range("J66").Value = 3420
$ws.Range("K66").Value = 11791.428
$ws.Range("L66").Value = 17100
$ws.Range("M66").Value = -8359.428
$ws.Range("N66").Value = -23964

$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41622

$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -128112

$ws.Range("H74").Value = 1074.1428
$ws.Range("I74").Value = 964.2778
$ws.Range("J74").Value = 1733.3334
$ws.Range("K74").Value = 964.2778
$ws.Range("L74").Value = 1733.3334
$ws.Range("M74").Value = -90.27779999999996
$ws.Range("N74").Value = -3481.3334

$ws.Range("H77").Value = 1074.1428
$ws.Range("I77").Value = 964.2778
$ws.Range("J77").Value = 1733.3334
$ws.Range("K77").Value = 4821.389
$ws.Range("L77").Value = 8666.666999999999
$ws.Range("M77").Value = -453.3890000000001
$ws.Range("N77").Value = -17402.667

$ws.Range("H136").Value = 3087.5
$ws.Range("I136").Value = 1042.5
$ws.Range("J136").Value = 8200
$ws.Range("K136").Value = 3127.5
$ws.Range("L136").Value = 24600
$ws.Range("M136").Value = -577.5
$ws.Range("N136").Value = -29700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 30780
$ws.Range("J51").Value = 30780
$ws.Range("L51").Value = 30780
$ws.Range("N51").Value = -31762

$ws.Range("H134").Value = 1997.7142
$ws.Range("I134").Value = 1429.7931
$ws.Range("J134").Value = 3264.6155
$ws.Range("K134").Value = 4289.379300000001
$ws.Range("L134").Value = 9793.8465
$ws.Range("M134").Value = -1754.379300000001
$ws.Range("N134").Value = -14863.8465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3129893.5
$ws.Range("I31").Value = 5557541.5
$ws.Range("J31").Value = 8632.071
$ws.Range("K31").Value = 5557541.5
$ws.Range("L31").Value = 8632.071
$ws.Range("M31").Value = -5557246.5
$ws.Range("N31").Value = -9222.071

$ws.Range("H34").Value = 3129893.5
$ws.Range("I34").Value = 5557541.5
$ws.Range("J34").Value = 8632.071
$ws.Range("K34").Value = 5557541.5
$ws.Range("L34").Value = 8632.071
$ws.Range("M34").Value = -5557339.5
$ws.Range("N34").Value = -9036.071

$ws.Range("H58").Value = 31254068
$ws.Range("I58").Value = 3550
$ws.Range("J58").Value = 41670908
$ws.Range("K58").Value = 3550
$ws.Range("L58").Value = 41670908
$ws.Range("M58").Value = -3347
$ws.Range("N58").Value = -41671314

$ws.Range("H132").Value = 3631.5186
$ws.Range("I132").Value = 2087.5386
$ws.Range("J132").Value = 5065.2144
$ws.Range("K132").Value = 6262.6158
$ws.Range("L132").Value = 15195.6432
$ws.Range("M132").Value = -3732.6158
$ws.Range("N132").Value = -20255.6432

$ws.Range("H134").Value = 2561.0715
$ws.Range("I134").Value = 1614.7646
$ws.Range("J134").Value = 4023.5454
$ws.Range("K134").Value = 4844.293799999999
$ws.Range("L134").Value = 12070.6362
$ws.Range("M134").Value = -2309.293799999999
$ws.Range("N134").Value = -17140.6362

$ws.Range("H136").Value = 31254068
$ws.Range("I136").Value = 3550
$ws.Range("J136").Value = 41670908
$ws.Range("K136").Value = 10650
$ws.Range("L136").Value = 125012724
$ws.Range("M136").Value = -8100
$ws.Range("N136").Value = -125017824

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 168.5
$ws.Range("J33").Value = 178.28572
$ws.Range("L33").Value = 1069.71432
$ws.Range("N33").Value = -1635.71432

$ws.Range("H120").Value = 16837.143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 37999.75
$ws.Range("J119").Value = 37999.75
$ws.Range("L119").Value = 37999.75
$ws.Range("N119").Value = -47675.75

$ws.Range("H126").Value = 4066.6667
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 4600
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 13800
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -18740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2307.25
$ws.Range("I7").Value = 1564.5
$ws.Range("J7").Value = 3050
$ws.Range("K7").Value = 1564.5
$ws.Range("L7").Value = 3050
$ws.Range("M7").Value = -1452.5
$ws.Range("N7").Value = -3274

$ws.Range("H126").Value = 2307.25
$ws.Range("I126").Value = 1564.5
$ws.Range("J126").Value = 3050
$ws.Range("K126").Value = 4693.5
$ws.Range("L126").Value = 9150
$ws.Range("M126").Value = -2223.5
$ws.Range("N126").Value = -14090

$ws.Range("H132").Value = 2338.9333
$ws.Range("I132").Value = 1716.0588
$ws.Range("J132").Value = 3153.4614
$ws.Range("K132").Value = 5148.1764
$ws.Range("L132").Value = 9460.3842
$ws.Range("M132").Value = -2618.1764
$ws.Range("N132").Value = -14520.3842

$ws.Range("H136").Value = 3335845
$ws.Range("I136").Value = 4547325
$ws.Range("K136").Value = 13641975
$ws.Range("M136").Value = -13639425

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 221310.94
$ws.Range("I132").Value = 323680.47
$ws.Range("J132").Value = 9747.200000000001
$ws.Range("K132").Value = 971041.4099999999
$ws.Range("L132").Value = 29241.6
$ws.Range("M132").Value = -968511.4099999999
$ws.Range("N132").Value = -34301.60000000001

$ws.Range("H136").Value = 1834.4517
$ws.Range("I136").Value = 879.53845
$ws.Range("K136").Value = 2638.61535
$ws.Range("M136").Value = -88.61535000000003
